$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 44-57, column E all currently show "A,4SET=1PC" -> fix typo to "A,4PC=1SET"
for ($r = 44; $r -le 57; $r++) {
    $ws.Cells.Item($r, 5).Value = "A,4PC=1SET"
}

# Set the frozen pane's scrolled top-left cell and the active selection within it
$ws.Activate()
$ws.Range("F52").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 4
